# Insert a new data row at row 167 (pushing existing rows 167:241 down to 168:242)
# and populate it with the new weekly price observation for Zanahoria
# (Vega Modelo de Temuco), matching the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(167).Insert()

$ws.Cells.Item(167, 1).Value  = 10
$ws.Cells.Item(167, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(167, 3).Value  = "La Araucanía"
$ws.Cells.Item(167, 4).Value  = 44553
$ws.Cells.Item(167, 5).Value  = 9
$ws.Cells.Item(167, 6).Value  = 100114013
$ws.Cells.Item(167, 7).Value  = "Zanahoria"
$ws.Cells.Item(167, 8).Value  = "Sin especificar"
$ws.Cells.Item(167, 9).Value  = "Primera"
$ws.Cells.Item(167, 10).Value = 305
$ws.Cells.Item(167, 11).Value = 6000
$ws.Cells.Item(167, 12).Value = 7000
$ws.Cells.Item(167, 13).Value = 6590
$ws.Cells.Item(167, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(167, 15).Value = "Región del Maule"
$ws.Cells.Item(167, 16).Value = 330
$ws.Cells.Item(167, 17).Value = 20
$ws.Cells.Item(167, 18).Value = "Hortaliza"
